$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1749.5
$ws.Range("I40").Value = 1325
$ws.Range("J40").Value = 1870.7858
$ws.Range("K40").Value = 1325
$ws.Range("L40").Value = 1870.7858
$ws.Range("M40").Value = -1150
$ws.Range("N40").Value = -2220.7858
$ws.Range("H74").Value = 3237.5
$ws.Range("J74").Value = 4500
$ws.Range("L74").Value = 4500
$ws.Range("N74").Value = -6372
$ws.Range("H77").Value = 3237.5
$ws.Range("J77").Value = 4500
$ws.Range("L77").Value = 22500
$ws.Range("N77").Value = -31860
$ws.Range("H80").Value = 814011.8
$ws.Range("I80").Value = 1516637.1
$ws.Range("K80").Value = 4549911.300000001
$ws.Range("M80").Value = -4548913.300000001
$ws.Range("H83").Value = 814011.8
$ws.Range("I83").Value = 1516637.1
$ws.Range("K83").Value = 13649733.9
$ws.Range("M83").Value = -13644741.9
$ws.Range("H132").Value = 2426.6494
$ws.Range("I132").Value = 2185.1858
$ws.Range("K132").Value = 6555.557400000001
$ws.Range("M132").Value = -4025.557400000001
$ws.Range("H137").Value = 25001868
$ws.Range("I137").Value = 55556896
$ws.Range("J137").Value = 2299.818
$ws.Range("K137").Value = 166670688
$ws.Range("L137").Value = 6899.454000000001
$ws.Range("M137").Value = -166668138
$ws.Range("N137").Value = -11999.454
$ws.Range("H138").Value = 1738.091
$ws.Range("I138").Value = 1418.3823
$ws.Range("K138").Value = 4255.1469
$ws.Range("M138").Value = 884.8531000000003
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2680.0454
$ws.Range("I32").Value = 2864.6667
$ws.Range("J32").Value = 1240
$ws.Range("K32").Value = 2864.6667
$ws.Range("L32").Value = 1240
$ws.Range("M32").Value = -2577.6667
$ws.Range("N32").Value = -1814
$ws.Range("H132").Value = 1554.7273
$ws.Range("I132").Value = 1542
$ws.Range("K132").Value = 4626
$ws.Range("M132").Value = -2096
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 458.91666
$ws.Range("I22").Value = 452.5
$ws.Range("J22").Value = 475.6
$ws.Range("K22").Value = 452.5
$ws.Range("L22").Value = 475.6
$ws.Range("M22").Value = -102.5
$ws.Range("N22").Value = -1175.6
$ws.Range("H31").Value = 4611.371
$ws.Range("I31").Value = 7337.2856
$ws.Range("K31").Value = 7337.2856
$ws.Range("M31").Value = -7042.2856
$ws.Range("H34").Value = 4611.371
$ws.Range("I34").Value = 7337.2856
$ws.Range("K34").Value = 7337.2856
$ws.Range("M34").Value = -7135.2856
$ws.Range("H58").Value = 1304.5957
$ws.Range("I58").Value = 1232.1212
$ws.Range("K58").Value = 1232.1212
$ws.Range("M58").Value = -1029.1212
$ws.Range("H132").Value = 3380.9167
$ws.Range("I132").Value = 2956.1052
$ws.Range("J132").Value = 4995.2
$ws.Range("K132").Value = 8868.3156
$ws.Range("L132").Value = 14985.6
$ws.Range("M132").Value = -6338.3156
$ws.Range("N132").Value = -20045.6
$ws.Range("H136").Value = 1304.5957
$ws.Range("I136").Value = 1232.1212
$ws.Range("K136").Value = 3696.3636
$ws.Range("M136").Value = -1146.3636
$ws.Range("H137").Value = 51444.445
$ws.Range("J137").Value = 52375
$ws.Range("L137").Value = 52375
$ws.Range("N137").Value = -62575
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 722.7778
$ws.Range("I17").Value = 162
$ws.Range("K17").Value = 486
$ws.Range("M17").Value = -317
$ws.Range("H23").Value = 57.666668
$ws.Range("I23").Value = 35.5
$ws.Range("K23").Value = 106.5
$ws.Range("M23").Value = 128.5
$ws.Range("H80").Value = 3025.625
$ws.Range("I80").Value = 1751.25
$ws.Range("K80").Value = 5253.75
$ws.Range("M80").Value = -4317.75
$ws.Range("H81").Value = 8630.666999999999
$ws.Range("J81").Value = 8630.666999999999
$ws.Range("L81").Value = 25892.001
$ws.Range("N81").Value = -28138.001
$ws.Range("H83").Value = 3025.625
$ws.Range("I83").Value = 1751.25
$ws.Range("K83").Value = 15761.25
$ws.Range("M83").Value = -11081.25
$ws.Range("H84").Value = 8630.666999999999
$ws.Range("J84").Value = 8630.666999999999
$ws.Range("L84").Value = 77676.003
$ws.Range("N84").Value = -88908.003
$ws.Range("H121").Value = 109756.73
$ws.Range("I121").Value = 17049
$ws.Range("J121").Value = 221006
$ws.Range("K121").Value = 51147
$ws.Range("L121").Value = 663018
$ws.Range("M121").Value = -49837
$ws.Range("N121").Value = -665638
$ws.Range("H122").Value = 990
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = ""
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4077.4736
$ws.Range("I80").Value = 3912.9167
$ws.Range("J80").Value = 4359.5713
$ws.Range("K80").Value = 3912.9167
$ws.Range("L80").Value = 4359.5713
$ws.Range("M80").Value = -2914.9167
$ws.Range("N80").Value = -6355.5713
$ws.Range("H83").Value = 4077.4736
$ws.Range("I83").Value = 3912.9167
$ws.Range("J83").Value = 4359.5713
$ws.Range("K83").Value = 19564.5835
$ws.Range("L83").Value = 21797.8565
$ws.Range("M83").Value = -14572.5835
$ws.Range("N83").Value = -31781.8565
$ws.Range("H102").Value = 2014.6666
$ws.Range("I102").Value = 1162.05
$ws.Range("J102").Value = 3719.9
$ws.Range("K102").Value = 1162.05
$ws.Range("L102").Value = 3719.9
$ws.Range("M102").Value = 459.95
$ws.Range("N102").Value = -6963.9
$ws.Range("H132").Value = 3368.7188
$ws.Range("J132").Value = 4913
$ws.Range("L132").Value = 14739
$ws.Range("N132").Value = -19799
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 698.4167
$ws.Range("I22").Value = 195.2
$ws.Range("K22").Value = 195.2
$ws.Range("M22").Value = 99.80000000000001
$ws.Range("H27").Value = 698.4167
$ws.Range("I27").Value = 195.2
$ws.Range("K27").Value = 195.2
$ws.Range("M27").Value = -88.19999999999999
$ws.Range("H46").Value = 2569.7097
$ws.Range("I46").Value = 1454
$ws.Range("K46").Value = 1454
$ws.Range("M46").Value = -1266
$ws.Range("H132").Value = 4127.452
$ws.Range("I132").Value = 2274.7097
$ws.Range("J132").Value = 9348.817999999999
$ws.Range("K132").Value = 6824.1291
$ws.Range("L132").Value = 28046.454
$ws.Range("M132").Value = -4294.1291
$ws.Range("N132").Value = -33106.454
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2436.5625
$ws.Range("I136").Value = 1138.64
$ws.Range("J136").Value = 7072
$ws.Range("K136").Value = 3415.92
$ws.Range("L136").Value = 21216
$ws.Range("M136").Value = -865.92
$ws.Range("N136").Value = -26316
